$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.780.68'
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").Value = '3.058.13'
$ws.Range("E3").Value = '  -0.25%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '555.66'
$ws.Range("E5").Value = '  +0.73%  '
$ws.Range("D6").Value = '142.73'
$ws.Range("E6").Value = '  +0.50%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '3.054.63'
$ws.Range("E8").Value = '  -0.20%  '
$ws.Range("D9").Value = '0.504'
$ws.Range("E9").Value = '  +0.30%  '
$ws.Range("D10").Value = '0.155'
$ws.Range("E10").Value = '  +2.17%  '
$ws.Range("D11").Value = '6.10'
$ws.Range("E11").Value = '  -6.35%  '
$ws.Range("D12").Value = '0.470'
$ws.Range("E12").Value = '  +2.68%  '
$ws.Range("D13").Value = '0.0000228'
$ws.Range("E13").Value = '  +0.29%  '
$ws.Range("D14").Value = '34.76'
$ws.Range("E14").Value = '  -0.13%  '
$ws.Range("D15").Value = '3.551.76'
$ws.Range("E15").Value = '  -0.40%  '
$ws.Range("D16").Value = '63.754.77'
$ws.Range("E16").Value = '  +0.52%  '
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").Value = '0.110'
$ws.Range("E17").Value = '  +0.65%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.050.55'
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("D19").Value = '6.69'
$ws.Range("E19").Value = '  -1.27%  '
$ws.Range("D20").Value = '478.13'
$ws.Range("E20").Value = '  -0.99%  '
$ws.Range("D21").Value = '13.96'
$ws.Range("E21").Value = '  +0.46%  '
$ws.Range("D22").Value = '0.674'
$ws.Range("E22").Value = '  -0.38%  '
$ws.Range("D23").Value = '7.50'
$ws.Range("E23").Value = '  +3.23%  '
$ws.Range("D24").Value = '14.04'
$ws.Range("E24").Value = '  +9.83%  '
$ws.Range("D25").Value = '80.96'
$ws.Range("E25").Value = '  +0.21%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("D27").Value = '2.78'
$ws.Range("E27").Value = '  +0.08%  '
$ws.Range("D28").Value = '7.94'
$ws.Range("E28").Value = '  +0.39%  '
$ws.Range("D29").Value = '2.04'
$ws.Range("E29").Value = '  +1.43%  '
$ws.Range("D30").Value = '0.998'
$ws.Range("E30").Value = '  -0.17%  '
$ws.Range("D31").Value = '26.12'
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("E32").Value = '  -2.35%  '
$ws.Range("D33").Value = '2.46'
$ws.Range("E33").Value = '  +0.54%  '
$ws.Range("D34").Value = '5.59'
$ws.Range("E34").Value = '  -1.87%  '
$ws.Range("D35").Value = '6.15'
$ws.Range("E35").Value = '  +2.61%  '
$ws.Range("D36").Value = '55.04'
$ws.Range("E36").Value = '  -0.56%  '
$ws.Range("D37").Value = '0.0405'
$ws.Range("E37").Value = '  +2.37%  '
$ws.Range("B38").Value = 'dogwifhat'
$ws.Range("C38").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D38").Value = '2.93'
$ws.Range("E38").Value = '  +13.99%  '
$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").Value = '439.40'
$ws.Range("E39").Value = '  -5.89%  '
$ws.Range("D40").Value = '0.0805'
$ws.Range("E40").Value = '  -2.15%  '
$ws.Range("D41").Value = '2.949.88'
$ws.Range("E41").Value = '  -1.98%  '
$ws.Range("D42").Value = '8.16'
$ws.Range("E42").Value = '  -0.47%  '
$ws.Range("E43").Value = '  -4.62%  '
$ws.Range("D44").Value = '28.26'
$ws.Range("E44").Value = '  +2.26%  '
$ws.Range("D45").Value = '0.258'
$ws.Range("E45").Value = '  +1.04%  '
$ws.Range("D47").Value = '2.12'
$ws.Range("E47").Value = '  +3.50%  '
$ws.Range("D48").Value = '0.112'
$ws.Range("E48").Value = '  +1.19%  '
$ws.Range("B49").Value = 'PEPE'
$ws.Range("C49").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D49").Value = '0.0₃0515'
$ws.Range("E49").Value = '  +1.26%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = '116.46'
$ws.Range("E50").Value = '  +0.10%  '
$ws.Range("D51").Value = '2.06'
$ws.Range("E51").Value = '  -0.74%  '
